# Besprechung in Stunden doku eingetragen
# Insert a new row for the "4. Teambesprechung (Spezifikation)" meeting
# above the existing row 7 (27.03.2017), shifting all following rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7 (pushes old row 7.. down to 8..)
$ws.Rows("7:7").Insert()
$ws.Rows("7:7").RowHeight = 15.75

# Fill in the new meeting entry
$ws.Range("A7").Value = 42821
$ws.Range("B7").Value = "16:00-17:00"
$ws.Range("C7").Value = "4. Teambesprechung (Spezifikation)"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "Andre, Bjorna, Johannes, Miel"

# Match the saved selection state from the edit
$ws.Range("F9").Select()
